$d = $word.ActiveDocument

# --- Hunk 1 -----------------------------------------------------------
# The blank paragraph right after "...CHANGE MY PROFILE PICTURE" (and
# before the "SEZIONE OPERAZIONI PAZIENTE" heading) gains an explicit
# en-US run language on its paragraph mark, alongside the existing
# italic formatting.
$r = $d.Content
$null = $r.Find.Execute("CHANGE MY PROFILE PICTURE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pPicture = $r.Paragraphs(1)
$pBlank = $pPicture.Next()
$pBlank.Range.LanguageID = "en-US"

# --- Hunk 2 -----------------------------------------------------------
# Locate the two "HAVE SETTINGS so that I can CHOOSE ..." bullets in
# "SEZIONE OPERAZIONI SEGRETARIO" (the paragraph objects get re-fetched
# via .Next() after each mutation below since inserting/deleting text
# can reseat them).
$r2 = $d.Content
$null = $r2.Find.Execute("CHOOSE MY CURRENT DOCTOR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pOffice = $r2.Paragraphs(1)
$pChoose = $pOffice.Next()

# Move the _GoBack bookmark from the start of $pOffice to the very end
# of $pChoose's text (right before its paragraph mark), then remove
# the now-superfluous empty paragraph that followed $pChoose.
#
# A zero-length Range placed exactly before a paragraph mark is not
# positioned reliably by Bookmarks.Add in this host, so a throwaway
# marker character is inserted, bookmarked, and then deleted again,
# leaving the (now collapsed) bookmark correctly anchored in place.
$insertPoint = $d.Range($pChoose.Range.End - 1, $pChoose.Range.End - 1)
$insertPoint.InsertAfter("X")

$pChoose = $pOffice.Next()
$markerRange = $d.Range($pChoose.Range.End - 2, $pChoose.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$pChoose = $pOffice.Next()
$markerRange = $d.Range($pChoose.Range.End - 2, $pChoose.Range.End - 1)
$markerRange.Delete()

$pChoose = $pOffice.Next()
$pEmpty = $pChoose.Next()
$pEmpty.Range.Delete()
